# Apply corrected IFRS financial figures for rows 2-6 (columns D:AJ)
# and remove the stale estimate columns from rows 7-9 (columns D:AI),
# matching the "error solve ifrs list" correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> @{ column number -> corrected value }
$newValues = @{
    2 = @{ 4 = 10467; 5 = 837; 6 = 837; 7 = 813; 8 = 623; 9 = 211; 10 = 412; 11 = 16341; 12 = 4715; 13 = 11625; 14 = 6146; 15 = 5479; 16 = 110; 17 = 1319; 18 = -1317; 19 = -263; 20 = 692; 21 = 627; 22 = 1647; 23 = 8; 24 = 5.96; 25 = 4.04; 26 = 3.81; 27 = 40.56; 28 = 9899.459999999999; 29 = 13710; 30 = 10.79; 31 = 323230; 32 = 0.46; 33 = 1500; 34 = 1.01; 35 = 13.5; 36 = 2191024 }
    3 = @{ 4 = 10519; 5 = 605; 6 = 605; 7 = 264; 8 = 82; 9 = 63; 10 = 19; 11 = 15926; 12 = 4309; 13 = 11617; 14 = 6164; 15 = 5453; 16 = 110; 17 = 836; 18 = -159; 19 = -646; 20 = 538; 21 = 299; 22 = 1072; 23 = 5.75; 24 = 0.78; 25 = 1.03; 26 = 0.51; 27 = 37.09; 28 = 9929.190000000001; 29 = 2893; 30 = 35.26; 31 = 324171; 32 = 0.31; 33 = 1500; 34 = 1.47; 35 = 45; 36 = 2191024 }
    4 = @{ 4 = 10953; 5 = 861; 6 = 861; 7 = 865; 8 = 676; 9 = 348; 10 = 328; 11 = 17065; 12 = 4814; 13 = 12251; 14 = 6497; 15 = 5754; 16 = 110; 17 = 1114; 18 = -1008; 19 = 427; 20 = 1050; 21 = 64; 22 = 1567; 23 = 7.86; 24 = 6.17; 25 = 5.5; 26 = 4.1; 27 = 39.3; 28 = 10244.71; 29 = 15890; 30 = 6.48; 31 = 341685; 32 = 0.3; 33 = 1500; 34 = 1.46; 35 = 8.19; 36 = 2191024 }
    5 = @{ 4 = 11978; 5 = 641; 6 = 641; 7 = 1045; 8 = 735; 9 = 437; 10 = 297; 11 = 18164; 12 = 5535; 13 = 12629; 14 = 6800; 15 = 5829; 16 = 110; 17 = 373; 18 = -505; 19 = 500; 20 = 1006; 21 = -634; 22 = 2176; 23 = 5.36; 24 = 6.13; 25 = 6.58; 26 = 4.17; 27 = 43.83; 28 = 10658.44; 29 = 19967; 30 = 5.73; 31 = 357595; 32 = 0.32; 33 = 1500; 34 = 1.31; 35 = 6.52; 36 = 2191024 }
    6 = @{ 4 = 16276; 5 = 1837; 6 = 1837; 7 = 1357; 8 = 1002; 9 = 501; 11 = 28403; 12 = 14564; 13 = 13838; 14 = 7364; 16 = 110; 17 = 2371; 18 = -3733; 19 = 1692; 20 = 1180; 21 = 1191; 22 = 8690; 23 = 11.28; 24 = 6.16; 25 = 7.07; 26 = 4.3; 27 = 105.25; 28 = 11190.53; 29 = 22849; 30 = 4.95; 31 = 387294; 32 = 0.29; 33 = 1750; 34 = 1.55; 35 = 6.65; 36 = 2191024 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $newValues[$row][$col]
    }
}

# row number -> list of column numbers whose cells must be emptied entirely
$clearCols = @{
    7 = @(4, 5, 7, 8, 9, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 23, 24, 25, 26, 27, 29, 30, 31, 32, 33, 34, 35)
    8 = @(4, 5, 7, 8, 9, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 23, 24, 25, 26, 27, 29, 30, 31, 32, 33, 34, 35)
    9 = @(4, 5, 7, 8, 9, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 23, 24, 25, 26, 27, 29, 30, 31, 32, 33, 34, 35)
}

foreach ($row in $clearCols.Keys) {
    foreach ($col in $clearCols[$row]) {
        $ws.Cells.Item($row, $col).ClearContents()
    }
}
